$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SD")

# Insert two new columns before column D (old D:K shifts to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Carry over number formatting/styles for the two new columns from column F
# (which now holds what used to be column D), so each row keeps its date /
# number / text style.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the two newly-inserted quarter columns (D = Q4'18, E = Q3'18)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 85100
$ws.Range("E8").Value = 97700
$ws.Range("D9").Value = 28500
$ws.Range("E9").Value = 29100
$ws.Range("D10").Value = 56600
$ws.Range("E10").Value = 68600
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 38000
$ws.Range("E15").Value = 36100
$ws.Range("D17").Value = 32300
$ws.Range("E17").Value = 85200
$ws.Range("D18").Value = 52800
$ws.Range("E18").Value = 12500
$ws.Range("D20").Value = 1400
$ws.Range("E20").Value = -800
$ws.Range("D21").Value = 92200
$ws.Range("E21").Value = 47800
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("D23").Value = 54200
$ws.Range("E23").Value = 11700
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 54200
$ws.Range("E26").Value = 11700
$ws.Range("D27").Value = 54200
$ws.Range("E27").Value = 11700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1400
$ws.Range("E32").Value = 800
$ws.Range("D33").Value = 54200
$ws.Range("E33").Value = 11700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 54200
$ws.Range("E35").Value = 11700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 17700
$ws.Range("E41").Value = 32600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 45500
$ws.Range("E43").Value = 54500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 10200
$ws.Range("E45").Value = 4600
$ws.Range("D46").Value = 73300
$ws.Range("E46").Value = 91600
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 949900
$ws.Range("E48").Value = 939500
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1100
$ws.Range("E52").Value = 1200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1024300
$ws.Range("E54").Value = 1032300
$ws.Range("D57").Value = 91000
$ws.Range("E57").Value = 97000
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 46200
$ws.Range("E59").Value = 92900
$ws.Range("D60").Value = 137200
$ws.Range("E60").Value = 189900
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 39400
$ws.Range("E62").Value = 49900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 176600
$ws.Range("E66").Value = 239800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -207500
$ws.Range("E72").Value = -261700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 847700
$ws.Range("E76").Value = 792500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 54200
$ws.Range("E81").Value = 11700
$ws.Range("D83").Value = 38000
$ws.Range("E83").Value = 36100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 36300
$ws.Range("E89").Value = 53100
$ws.Range("D91").Value = 122100
$ws.Range("E91").Value = -51500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -51100
$ws.Range("E94").Value = -50600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -14800
$ws.Range("E102").Value = 2500

# A few historical "Capital Expenditures" figures were restated for three of
# the existing quarters when this row was re-derived.
$ws.Range("H91").Value = -100
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = -4600
